# Auto-generated edit script: apply scheduled market-price/profit updates
# to the Leve profit tables across all 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 2454.75  # H19: 2454.7917 -> 2454.75
$ws.Cells.Item(19, 9).Value = 2151.75  # I19: 2120.7693 -> 2151.75
$ws.Cells.Item(19, 10).Value = 2757.75  # J19: 2849.5454 -> 2757.75
$ws.Cells.Item(19, 11).Value = 2151.75  # K19: 2120.7693 -> 2151.75
$ws.Cells.Item(19, 12).Value = 2757.75  # L19: 2849.5454 -> 2757.75
$ws.Cells.Item(19, 13).Value = -1976.75  # M19: -1945.7693 -> -1976.75
$ws.Cells.Item(19, 14).Value = -3107.75  # N19: -3199.5454 -> -3107.75

$ws.Cells.Item(39, 8).Value = 3346.1875  # H39: 3544.6 -> 3346.1875
$ws.Cells.Item(39, 10).Value = 6256  # J39: 7096.857 -> 6256
$ws.Cells.Item(39, 12).Value = 18768  # L39: 21290.571 -> 18768
$ws.Cells.Item(39, 14).Value = -19360  # N39: -21882.571 -> -19360

$ws.Cells.Item(86, 8).Value = 3640  # H86: 3818.4167 -> 3640
$ws.Cells.Item(86, 9).Value = 3931.7  # I86: 4202 -> 3931.7
$ws.Cells.Item(86, 11).Value = 3931.7  # K86: 4202 -> 3931.7
$ws.Cells.Item(86, 13).Value = -2808.7  # M86: -3079 -> -2808.7

$ws.Cells.Item(89, 8).Value = 3640  # H89: 3818.4167 -> 3640
$ws.Cells.Item(89, 9).Value = 3931.7  # I89: 4202 -> 3931.7
$ws.Cells.Item(89, 11).Value = 19658.5  # K89: 21010 -> 19658.5
$ws.Cells.Item(89, 13).Value = -14042.5  # M89: -15394 -> -14042.5

$ws.Cells.Item(101, 8).Value = 487.54544  # H101: 471.91666 -> 487.54544
$ws.Cells.Item(101, 9).Value = 502  # I101: 468.33334 -> 502
$ws.Cells.Item(101, 11).Value = 1506  # K101: 1405.00002 -> 1506
$ws.Cells.Item(101, 13).Value = 116  # M101: 216.9999800000001 -> 116

$ws.Cells.Item(106, 8).Value = 9557.143  # H106: 9833.333000000001 -> 9557.143
$ws.Cells.Item(106, 9).Value = 8450  # I106: 9000 -> 8450
$ws.Cells.Item(106, 11).Value = 8450  # K106: 9000 -> 8450
$ws.Cells.Item(106, 13).Value = -7819  # M106: -8369 -> -7819

$ws.Cells.Item(113, 8).Value = 2559.8  # H113: 2236.25 -> 2559.8
$ws.Cells.Item(113, 9).Value = 2559.8  # I113: 2383.1667 -> 2559.8
$ws.Cells.Item(113, 10).Value = 0  # J113: 1795.5 -> 0
$ws.Cells.Item(113, 11).Value = 2559.8  # K113: 2383.1667 -> 2559.8
$ws.Cells.Item(113, 12).Value = 0  # L113: 1795.5 -> 0
$ws.Cells.Item(113, 13).Value = 694.1999999999998  # M113: 870.8332999999998 -> 694.1999999999998
$ws.Cells.Item(113, 14).ClearContents()  # N113: -8303.5 -> (removed)

$ws.Cells.Item(135, 8).Value = 18907.857  # H135: 21651 -> 18907.857
$ws.Cells.Item(135, 9).Value = 5071.2  # I135: 5726.75 -> 5071.2
$ws.Cells.Item(135, 11).Value = 45640.8  # K135: 51540.75 -> 45640.8
$ws.Cells.Item(135, 13).Value = -43105.8  # M135: -49005.75 -> -43105.8

$ws.Cells.Item(141, 8).Value = 6131.091  # H141: 6146 -> 6131.091
$ws.Cells.Item(141, 9).Value = 6051  # I141: 6059.25 -> 6051
$ws.Cells.Item(141, 10).Value = 6491.5  # J141: 6493 -> 6491.5
$ws.Cells.Item(141, 11).Value = 18153  # K141: 18177.75 -> 18153
$ws.Cells.Item(141, 12).Value = 19474.5  # L141: 19479 -> 19474.5
$ws.Cells.Item(141, 13).Value = -12973  # M141: -12997.75 -> -12973
$ws.Cells.Item(141, 14).Value = -29834.5  # N141: -29839 -> -29834.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 10767.541  # H32: 11365.143 -> 10767.541
$ws.Cells.Item(32, 9).Value = 1015.53845  # I32: 1074.375 -> 1015.53845
$ws.Cells.Item(32, 11).Value = 1015.53845  # K32: 1074.375 -> 1015.53845
$ws.Cells.Item(32, 13).Value = -728.53845  # M32: -787.375 -> -728.53845

$ws.Cells.Item(61, 8).Value = 6636.394  # H61: 6950.567 -> 6636.394
$ws.Cells.Item(61, 9).Value = 5307.2607  # I61: 5579.15 -> 5307.2607
$ws.Cells.Item(61, 11).Value = 5307.2607  # K61: 5579.15 -> 5307.2607
$ws.Cells.Item(61, 13).Value = -5095.2607  # M61: -5367.15 -> -5095.2607

$ws.Cells.Item(63, 8).Value = 3125.9048  # H63: 3318.5715 -> 3125.9048
$ws.Cells.Item(63, 9).Value = 3125.9048  # I63: 3318.5715 -> 3125.9048
$ws.Cells.Item(63, 11).Value = 3125.9048  # K63: 3318.5715 -> 3125.9048
$ws.Cells.Item(63, 13).Value = -2439.9048  # M63: -2632.5715 -> -2439.9048

$ws.Cells.Item(66, 8).Value = 3125.9048  # H66: 3318.5715 -> 3125.9048
$ws.Cells.Item(66, 9).Value = 3125.9048  # I66: 3318.5715 -> 3125.9048
$ws.Cells.Item(66, 11).Value = 15629.524  # K66: 16592.8575 -> 15629.524
$ws.Cells.Item(66, 13).Value = -12197.524  # M66: -13160.8575 -> -12197.524

$ws.Cells.Item(74, 8).Value = 2597.1765  # H74: 2597.2942 -> 2597.1765
$ws.Cells.Item(74, 9).Value = 1850.5834  # I74: 1850.75 -> 1850.5834
$ws.Cells.Item(74, 11).Value = 1850.5834  # K74: 1850.75 -> 1850.5834
$ws.Cells.Item(74, 13).Value = -976.5834  # M74: -976.75 -> -976.5834

$ws.Cells.Item(77, 8).Value = 2597.1765  # H77: 2597.2942 -> 2597.1765
$ws.Cells.Item(77, 9).Value = 1850.5834  # I77: 1850.75 -> 1850.5834
$ws.Cells.Item(77, 11).Value = 9252.916999999999  # K77: 9253.75 -> 9252.916999999999
$ws.Cells.Item(77, 13).Value = -4884.916999999999  # M77: -4885.75 -> -4884.916999999999

$ws.Cells.Item(132, 8).Value = 2676.16  # H132: 2634.6924 -> 2676.16
$ws.Cells.Item(132, 9).Value = 2682.3  # I132: 2630.6667 -> 2682.3
$ws.Cells.Item(132, 11).Value = 8046.900000000001  # K132: 7892.000100000001 -> 8046.900000000001
$ws.Cells.Item(132, 13).Value = -5516.900000000001  # M132: -5362.000100000001 -> -5516.900000000001

$ws.Cells.Item(136, 8).Value = 6636.394  # H136: 6950.567 -> 6636.394
$ws.Cells.Item(136, 9).Value = 5307.2607  # I136: 5579.15 -> 5307.2607
$ws.Cells.Item(136, 11).Value = 15921.7821  # K136: 16737.45 -> 15921.7821
$ws.Cells.Item(136, 13).Value = -13371.7821  # M136: -14187.45 -> -13371.7821

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2859.6667  # H86: 2598.4443 -> 2859.6667
$ws.Cells.Item(86, 9).Value = 3164.75  # I86: 2836.1667 -> 3164.75
$ws.Cells.Item(86, 10).Value = 2249.5  # J86: 2123 -> 2249.5
$ws.Cells.Item(86, 11).Value = 3164.75  # K86: 2836.1667 -> 3164.75
$ws.Cells.Item(86, 12).Value = 2249.5  # L86: 2123 -> 2249.5
$ws.Cells.Item(86, 13).Value = -2041.75  # M86: -1713.1667 -> -2041.75
$ws.Cells.Item(86, 14).Value = -4495.5  # N86: -4369 -> -4495.5

$ws.Cells.Item(89, 8).Value = 2859.6667  # H89: 2598.4443 -> 2859.6667
$ws.Cells.Item(89, 9).Value = 3164.75  # I89: 2836.1667 -> 3164.75
$ws.Cells.Item(89, 10).Value = 2249.5  # J89: 2123 -> 2249.5
$ws.Cells.Item(89, 11).Value = 15823.75  # K89: 14180.8335 -> 15823.75
$ws.Cells.Item(89, 12).Value = 11247.5  # L89: 10615 -> 11247.5
$ws.Cells.Item(89, 13).Value = -10207.75  # M89: -8564.833500000001 -> -10207.75
$ws.Cells.Item(89, 14).Value = -22479.5  # N89: -21847 -> -22479.5

$ws.Cells.Item(125, 8).Value = 88332.336  # H125: 88332.5 -> 88332.336
$ws.Cells.Item(125, 10).Value = 88332.336  # J125: 88332.5 -> 88332.336
$ws.Cells.Item(125, 12).Value = 88332.336  # L125: 88332.5 -> 88332.336
$ws.Cells.Item(125, 14).Value = -98172.336  # N125: -98172.5 -> -98172.336

$ws.Cells.Item(134, 8).Value = 2082.8333  # H134: 2030.3478 -> 2082.8333
$ws.Cells.Item(134, 9).Value = 1863.85  # I134: 1834.35 -> 1863.85
$ws.Cells.Item(134, 10).Value = 3177.75  # J134: 3337 -> 3177.75
$ws.Cells.Item(134, 11).Value = 5591.549999999999  # K134: 5503.049999999999 -> 5591.549999999999
$ws.Cells.Item(134, 12).Value = 9533.25  # L134: 10011 -> 9533.25
$ws.Cells.Item(134, 13).Value = -3056.549999999999  # M134: -2968.049999999999 -> -3056.549999999999
$ws.Cells.Item(134, 14).Value = -14603.25  # N134: -15081 -> -14603.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 40791.473  # H16: 40791.527 -> 40791.473
$ws.Cells.Item(16, 9).Value = 30638.166  # I16: 30638.25 -> 30638.166
$ws.Cells.Item(16, 11).Value = 30638.166  # K16: 30638.25 -> 30638.166
$ws.Cells.Item(16, 13).Value = -30351.166  # M16: -30351.25 -> -30351.166

$ws.Cells.Item(31, 8).Value = 6497.5757  # H31: 6997.567 -> 6497.5757
$ws.Cells.Item(31, 9).Value = 6733.6294  # I31: 7388.125 -> 6733.6294
$ws.Cells.Item(31, 11).Value = 6733.6294  # K31: 7388.125 -> 6733.6294
$ws.Cells.Item(31, 13).Value = -6438.6294  # M31: -7093.125 -> -6438.6294

$ws.Cells.Item(34, 8).Value = 6497.5757  # H34: 6997.567 -> 6497.5757
$ws.Cells.Item(34, 9).Value = 6733.6294  # I34: 7388.125 -> 6733.6294
$ws.Cells.Item(34, 11).Value = 6733.6294  # K34: 7388.125 -> 6733.6294
$ws.Cells.Item(34, 13).Value = -6531.6294  # M34: -7186.125 -> -6531.6294

$ws.Cells.Item(107, 8).Value = 1467.5  # H107: 621.4 -> 1467.5
$ws.Cells.Item(107, 9).Value = 1962.7  # I107: 676.2353000000001 -> 1962.7
$ws.Cells.Item(107, 10).Value = 642.1667  # J107: 504.875 -> 642.1667
$ws.Cells.Item(107, 11).Value = 1962.7  # K107: 676.2353000000001 -> 1962.7
$ws.Cells.Item(107, 12).Value = 642.1667  # L107: 504.875 -> 642.1667
$ws.Cells.Item(107, 13).Value = -42.70000000000005  # M107: 1243.7647 -> -42.70000000000005
$ws.Cells.Item(107, 14).Value = -4482.1667  # N107: -4344.875 -> -4482.1667

$ws.Cells.Item(113, 8).Value = 40791.473  # H113: 40791.527 -> 40791.473
$ws.Cells.Item(113, 9).Value = 30638.166  # I113: 30638.25 -> 30638.166
$ws.Cells.Item(113, 11).Value = 30638.166  # K113: 30638.25 -> 30638.166
$ws.Cells.Item(113, 13).Value = -28468.166  # M113: -28468.25 -> -28468.166

$ws.Cells.Item(125, 8).Value = 67223  # H125: 76532 -> 67223
$ws.Cells.Item(125, 9).Value = 40000  # I125: 0 -> 40000
$ws.Cells.Item(125, 10).Value = 76297.336  # J125: 76532 -> 76297.336
$ws.Cells.Item(125, 11).Value = 40000  # K125: 0 -> 40000
$ws.Cells.Item(125, 12).Value = 76297.336  # L125: 76532 -> 76297.336
$ws.Cells.Item(125, 13).Value = -37540  # M125: None -> -37540
$ws.Cells.Item(125, 14).Value = -81217.336  # N125: -81452 -> -81217.336

$ws.Cells.Item(132, 8).Value = 5970.88  # H132: 5837.346 -> 5970.88
$ws.Cells.Item(132, 9).Value = 4781.294  # I132: 4654.5 -> 4781.294
$ws.Cells.Item(132, 11).Value = 14343.882  # K132: 13963.5 -> 14343.882
$ws.Cells.Item(132, 13).Value = -11813.882  # M132: -11433.5 -> -11813.882

$ws.Cells.Item(139, 8).Value = 70690  # H139: 70990 -> 70690
$ws.Cells.Item(139, 10).Value = 84365.164  # J139: 84865.164 -> 84365.164
$ws.Cells.Item(139, 12).Value = 84365.164  # L139: 84865.164 -> 84365.164
$ws.Cells.Item(139, 14).Value = -94645.164  # N139: -95145.164 -> -94645.164

$ws.Cells.Item(141, 8).Value = 89424.75  # H141: 86339.60000000001 -> 89424.75
$ws.Cells.Item(141, 10).Value = 89424.75  # J141: 86339.60000000001 -> 89424.75
$ws.Cells.Item(141, 12).Value = 89424.75  # L141: 86339.60000000001 -> 89424.75
$ws.Cells.Item(141, 14).Value = -99784.75  # N141: -96699.60000000001 -> -99784.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(18, 8).Value = 266.42856  # H18: 255.875 -> 266.42856
$ws.Cells.Item(18, 9).Value = 266.42856  # I18: 268.2857 -> 266.42856
$ws.Cells.Item(18, 10).Value = 0  # J18: 169 -> 0
$ws.Cells.Item(18, 11).Value = 799.28568  # K18: 804.8571000000001 -> 799.28568
$ws.Cells.Item(18, 12).Value = 0  # L18: 507 -> 0
$ws.Cells.Item(18, 13).Value = -630.28568  # M18: -635.8571000000001 -> -630.28568
$ws.Cells.Item(18, 14).ClearContents()  # N18: -845 -> (removed)

$ws.Cells.Item(50, 8).Value = 611.5294  # H50: 533.0625 -> 611.5294
$ws.Cells.Item(50, 9).Value = 347.9  # I50: 326.18182 -> 347.9
$ws.Cells.Item(50, 10).Value = 988.1429000000001  # J50: 988.2 -> 988.1429000000001
$ws.Cells.Item(50, 11).Value = 1043.7  # K50: 978.54546 -> 1043.7
$ws.Cells.Item(50, 12).Value = 2964.4287  # L50: 2964.6 -> 2964.4287
$ws.Cells.Item(50, 13).Value = -562.6999999999998  # M50: -497.54546 -> -562.6999999999998
$ws.Cells.Item(50, 14).Value = -3926.4287  # N50: -3926.6 -> -3926.4287

$ws.Cells.Item(53, 8).Value = 611.5294  # H53: 533.0625 -> 611.5294
$ws.Cells.Item(53, 9).Value = 347.9  # I53: 326.18182 -> 347.9
$ws.Cells.Item(53, 10).Value = 988.1429000000001  # J53: 988.2 -> 988.1429000000001
$ws.Cells.Item(53, 11).Value = 1043.7  # K53: 978.54546 -> 1043.7
$ws.Cells.Item(53, 12).Value = 2964.4287  # L53: 2964.6 -> 2964.4287
$ws.Cells.Item(53, 13).Value = -562.6999999999998  # M53: -497.54546 -> -562.6999999999998
$ws.Cells.Item(53, 14).Value = -3926.4287  # N53: -3926.6 -> -3926.4287

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(43, 8).Value = 6174.1113  # H43: 5210.636 -> 6174.1113
$ws.Cells.Item(43, 9).Value = 6174.1113  # I43: 5210.636 -> 6174.1113
$ws.Cells.Item(43, 11).Value = 6174.1113  # K43: 5210.636 -> 6174.1113
$ws.Cells.Item(43, 13).Value = -6023.1113  # M43: -5059.636 -> -6023.1113

$ws.Cells.Item(70, 8).Value = 5765.636  # H70: 5907.364 -> 5765.636
$ws.Cells.Item(70, 9).Value = 5746.143  # I70: 5968.857 -> 5746.143
$ws.Cells.Item(70, 11).Value = 5746.143  # K70: 5968.857 -> 5746.143
$ws.Cells.Item(70, 13).Value = -5476.143  # M70: -5698.857 -> -5476.143

$ws.Cells.Item(73, 8).Value = 5765.636  # H73: 5907.364 -> 5765.636
$ws.Cells.Item(73, 9).Value = 5746.143  # I73: 5968.857 -> 5746.143
$ws.Cells.Item(73, 11).Value = 5746.143  # K73: 5968.857 -> 5746.143
$ws.Cells.Item(73, 13).Value = -4810.143  # M73: -5032.857 -> -4810.143

$ws.Cells.Item(102, 8).Value = 2281.1177  # H102: 2285.5 -> 2281.1177
$ws.Cells.Item(102, 9).Value = 2341.5715  # I102: 2326.3572 -> 2341.5715
$ws.Cells.Item(102, 10).Value = 1999  # J102: 1999.5 -> 1999
$ws.Cells.Item(102, 11).Value = 2341.5715  # K102: 2326.3572 -> 2341.5715
$ws.Cells.Item(102, 12).Value = 1999  # L102: 1999.5 -> 1999
$ws.Cells.Item(102, 13).Value = -719.5715  # M102: -704.3571999999999 -> -719.5715
$ws.Cells.Item(102, 14).Value = -5243  # N102: -5243.5 -> -5243

$ws.Cells.Item(126, 8).Value = 3643.7334  # H126: 3572.1875 -> 3643.7334
$ws.Cells.Item(126, 9).Value = 2459.2  # I126: 2462.818 -> 2459.2
$ws.Cells.Item(126, 11).Value = 7377.599999999999  # K126: 7388.454000000001 -> 7377.599999999999
$ws.Cells.Item(126, 13).Value = -4907.599999999999  # M126: -4918.454000000001 -> -4907.599999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 6676.7  # H46: 6797 -> 6676.7
$ws.Cells.Item(46, 9).Value = 6363  # I46: 6526.7 -> 6363
$ws.Cells.Item(46, 11).Value = 6363  # K46: 6526.7 -> 6363
$ws.Cells.Item(46, 13).Value = -6175  # M46: -6338.7 -> -6175

$ws.Cells.Item(68, 8).Value = 3650.2222  # H68: 3385.2 -> 3650.2222
$ws.Cells.Item(68, 9).Value = 4064.5715  # I68: 3681.5 -> 4064.5715
$ws.Cells.Item(68, 11).Value = 4064.5715  # K68: 3681.5 -> 4064.5715
$ws.Cells.Item(68, 13).Value = -3315.5715  # M68: -2932.5 -> -3315.5715

$ws.Cells.Item(71, 8).Value = 3650.2222  # H71: 3385.2 -> 3650.2222
$ws.Cells.Item(71, 9).Value = 4064.5715  # I71: 3681.5 -> 4064.5715
$ws.Cells.Item(71, 11).Value = 20322.8575  # K71: 18407.5 -> 20322.8575
$ws.Cells.Item(71, 13).Value = -16578.8575  # M71: -14663.5 -> -16578.8575

$ws.Cells.Item(114, 8).Value = 69999.5  # H114: 70000 -> 69999.5
$ws.Cells.Item(114, 10).Value = 69999.5  # J114: 70000 -> 69999.5
$ws.Cells.Item(114, 12).Value = 69999.5  # L114: 70000 -> 69999.5
$ws.Cells.Item(114, 14).Value = -78677.5  # N114: -78678 -> -78677.5

$ws.Cells.Item(122, 8).Value = 6172.4546  # H122: 6477.6665 -> 6172.4546
$ws.Cells.Item(122, 9).Value = 5150.1665  # I122: 5325.75 -> 5150.1665
$ws.Cells.Item(122, 11).Value = 15450.4995  # K122: 15977.25 -> 15450.4995
$ws.Cells.Item(122, 13).Value = -13000.4995  # M122: -13527.25 -> -13000.4995

$ws.Cells.Item(125, 8).Value = 89123.5  # H125: 89373.625 -> 89123.5
$ws.Cells.Item(125, 10).Value = 89123.5  # J125: 89373.625 -> 89123.5
$ws.Cells.Item(125, 12).Value = 89123.5  # L125: 89373.625 -> 89123.5
$ws.Cells.Item(125, 14).Value = -98963.5  # N125: -99213.625 -> -98963.5

$ws.Cells.Item(132, 8).Value = 40488.844  # H132: 42627.223 -> 40488.844
$ws.Cells.Item(132, 10).Value = 4378.75  # J132: 5172.3335 -> 4378.75
$ws.Cells.Item(132, 12).Value = 13136.25  # L132: 15517.0005 -> 13136.25
$ws.Cells.Item(132, 14).Value = -18196.25  # N132: -20577.0005 -> -18196.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 473.43478  # H113: 509.57144 -> 473.43478
$ws.Cells.Item(113, 10).Value = 302.16666  # J113: 406.25 -> 302.16666
$ws.Cells.Item(113, 12).Value = 906.4999799999999  # L113: 1218.75 -> 906.4999799999999
$ws.Cells.Item(113, 14).Value = -5246.49998  # N113: -5558.75 -> -5246.49998

$ws.Cells.Item(136, 8).Value = 4782.6772  # H136: 6564.0713 -> 4782.6772
$ws.Cells.Item(136, 9).Value = 4883.815  # I136: 5538.8696 -> 4883.815
$ws.Cells.Item(136, 10).Value = 4100  # J136: 11280 -> 4100
$ws.Cells.Item(136, 11).Value = 14651.445  # K136: 16616.6088 -> 14651.445
$ws.Cells.Item(136, 12).Value = 12300  # L136: 33840 -> 12300
$ws.Cells.Item(136, 13).Value = -12101.445  # M136: -14066.6088 -> -12101.445
$ws.Cells.Item(136, 14).Value = -17400  # N136: -38940 -> -17400
